# Updates cryptos list cell values to match the latest scrape (GitHub Actions run).
# Source data is plain text (prices keep their "thousand-dot" display format, e.g.
# "58.771.14", and are NOT real decimal numbers), so every target cell is written as
# text. Cells whose new text would otherwise be auto-recognised by Excel as a number
# (e.g. "0.999", "1.00") are first force-formatted as Text ("@") so COM does not
# silently coerce them into numeric values.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ref -> new value, in row order matching the source diff
$updates = [ordered]@{
    'D2' = '58.771.14'
    'E2' = '  +0.57%  '
    'D3' = '3.154.65'
    'E3' = '  +1.69%  '
    'D4' = '0.999'
    'E4' = '  -0.06%  '
    'D5' = '535.78'
    'E5' = '  +1.60%  '
    'D6' = '143.93'
    'E6' = '  +1.22%  '
    'E7' = '  -0.10%  '
    'D8' = '3.154.93'
    'E8' = '  +1.73%  '
    'D9' = '0.449'
    'E9' = '  +2.36%  '
    'D10' = '7.18'
    'E10' = '  -1.73%  '
    'E11' = '  +1.07%  '
    'D12' = '0.395'
    'E12' = '  +3.03%  '
    'D13' = '3.692.26'
    'E13' = '  +1.56%  '
    'E14' = '  +3.11%  '
    'D15' = '25.89'
    'E15' = '  -3.43%  '
    'E16' = '  +1.23%  '
    'D17' = '58.683.47'
    'E17' = '  +0.32%  '
    'D18' = '3.142.03'
    'E18' = '  +1.06%  '
    'E19' = '  +0.94%  '
    'E20' = '  +0.56%  '
    'E21' = '  -0.37%  '
    'D22' = '344.73'
    'E22' = '  +1.24%  '
    'E23' = '  +0.17%  '
    'D24' = '0.516'
    'E24' = '  +2.68%  '
    'E25' = '  +2.94%  '
    'D26' = '0.171'
    'E26' = '  +0.45%  '
    'D27' = '1.00'
    'E27' = '  +0.17%  '
    'D28' = '0.0₃0945'
    'E28' = '  +3.74%  '
    'D29' = '7.57'
    'E29' = '  +3.64%  '
    'E30' = '  -0.68%  '
    'E31' = '  +0.07%  '
    'E32' = '  +2.22%  '
    'B33' = 'EthereumClassic'
    'C33' = 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc'
    'D33' = '21.29'
    'E33' = '  +1.71%  '
    'B34' = 'Fetch.AI'
    'C34' = 'https://coinranking.com/coin/AWma-WzFHmKVQ+fetchai-fet'
    'D34' = '1.21'
    'E34' = '  +0.83%  '
    'D35' = '4.85'
    'E35' = '  +5.13%  '
    'D36' = '158.64'
    'E36' = '  +2.67%  '
    'D37' = '6.28'
    'E37' = '  +3.97%  '
    'D38' = '26.42'
    'E38' = '  -2.52%  '
    'E39' = '  -1.22%  '
    'D40' = '1.67'
    'E40' = '  +12.54%  '
    'E41' = '  +0.10%  '
    'D42' = '0.713'
    'E42' = '  +5.03%  '
    'D43' = '4.05'
    'E43' = '  +4.32%  '
    'D44' = '3.190.77'
    'E44' = '  +1.50%  '
    'D45' = '36.82'
    'E45' = '  +0.01%  '
    'D46' = '0.999'
    'E46' = '  -0.03%  '
    'D47' = '2.325.65'
    'E47' = '  +1.68%  '
    'E48' = '  +3.90%  '
    'E49' = '  +6.24%  '
    'D50' = '20.85'
    'E50' = '  +0.61%  '
    'D51' = '6.09'
    'E51' = '  +1.88%  '
}

foreach ($ref in $updates.Keys) {
    $value = $updates[$ref]
    $range = $ws.Range($ref)

    # Numeric-looking text (plain decimals like "0.999" or "1.00") would be
    # auto-converted to a real number by the Value setter, losing the trailing
    # zero / fixed formatting. Force the cell to Text first to keep it literal,
    # matching the inline-string cells in the source workbook.
    $isNumericLooking = $value -match '^\s*[+-]?\d+(\.\d+)?\s*$'
    if ($isNumericLooking) {
        $range.NumberFormat = '@'
    }

    $range.Value = $value
}
